$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.055.09"
$ws.Range("E2").Value = "  -1.41%  "

$ws.Range("D3").Value = "2.472.79"
$ws.Range("E3").Value = "  -2.57%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.50%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.88%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  -2.03%  "

$ws.Range("D9").Value = "2.472.58"
$ws.Range("E9").Value = "  -2.54%  "

$ws.Range("E10").Value = "  -2.68%  "

$ws.Range("E11").Value = "  -0.97%  "

$ws.Range("E12").Value = "  -2.52%  "

$ws.Range("E13").Value = "  -3.84%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.24%  "

$ws.Range("E15").Value = "  -3.09%  "

$ws.Range("D16").Value = "66.993.65"
$ws.Range("E16").Value = "  -1.48%  "

$ws.Range("E17").Value = "  -4.25%  "

$ws.Range("D18").Value = "2.477.52"
$ws.Range("E18").Value = "  -0.95%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.80%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "351.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.09%  "

$ws.Range("E22").Value = "  -2.72%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.95%  "

$ws.Range("E26").Value = "  -4.83%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.39%  "

$ws.Range("E28").Value = "  -57.79%  "

$ws.Range("D29").Value = "2.597.40"
$ws.Range("E29").Value = "  -2.60%  "

$ws.Range("E30").Value = "  -6.74%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "509.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.27%  "

$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.60%  "

$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.46%  "

$ws.Range("E35").Value = "  +0.08%  "

$ws.Range("E36").Value = "  -0.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.116"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -10.18%  "

$ws.Range("E38").Value = "  +0.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.32"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.55%  "

$ws.Range("E40").Value = "  -7.60%  "

$ws.Range("E41").Value = "  -5.17%  "

$ws.Range("E42").Value = "  +0.05%  "

$ws.Range("E43").Value = "  -6.24%  "

$ws.Range("E44").Value = "  -6.54%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "140.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.77%  "

$ws.Range("E48").Value = "  -6.62%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.97%  "

$ws.Range("E50").Value = "  -10.70%  "

$ws.Range("E51").Value = "  -7.27%  "
